$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.05619466666666667
$ws.Cells.Item(2, 9).Value = 0.04986276087265156
$ws.Cells.Item(2, 10).Value = 0.07297477932340853
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.874784666666667
$ws.Cells.Item(2, 14).Value = 8.624354
$ws.Cells.Item(2, 15).Value = 0.1187109652550681
$ws.Cells.Item(2, 16).Value = 0.121184727686443
$ws.Cells.Item(2, 17).Value = 0.1615475660817778
$ws.Cells.Item(2, 18).Value = 1.453928094736
$ws.Cells.Item(2, 19).Value = 0.005919256473475107
$ws.Cells.Item(2, 20).Value = 0.008843428760285535
$ws.Cells.Item(3, 7).Value = 0.05619466666666667
$ws.Cells.Item(3, 9).Value = 0.04986276087265156
$ws.Cells.Item(3, 10).Value = 0.07297477932340853
$ws.Cells.Item(3, 15).Value = 0.4442422727481699
$ws.Cells.Item(3, 16).Value = 0.4534996302499962
$ws.Cells.Item(3, 17).Value = 0.6045461576266666
$ws.Cells.Item(3, 18).Value = 5.44091541864
$ws.Cells.Item(3, 19).Value = 0.02215114621556525
$ws.Cells.Item(3, 20).Value = 0.03309403544074083
$ws.Cells.Item(4, 7).Value = 0.05619466666666667
$ws.Cells.Item(4, 9).Value = 0.04986276087265156
$ws.Cells.Item(4, 10).Value = 0.07297477932340853
$ws.Cells.Item(4, 13).Value = 5.147441999999999
$ws.Cells.Item(4, 14).Value = 15.442326
$ws.Cells.Item(4, 15).Value = 0.2125577666737049
$ws.Cells.Item(4, 16).Value = 0.2169871588243338
$ws.Cells.Item(4, 17).Value = 0.289258787376
$ws.Cells.Item(4, 18).Value = 2.603329086384
$ws.Cells.Item(4, 19).Value = 0.01059871709127581
$ws.Cells.Item(4, 20).Value = 0.01583459003121916
$ws.Cells.Item(5, 7).Value = 0.05619466666666667
$ws.Cells.Item(5, 9).Value = 0.04986276087265156
$ws.Cells.Item(5, 10).Value = 0.07297477932340853
$ws.Cells.Item(5, 13).Value = 1.483016
$ws.Cells.Item(5, 14).Value = 2.966032
$ws.Cells.Item(5, 15).Value = 0.06123946008548931
$ws.Cells.Item(5, 16).Value = 0.04167706708575228
$ws.Cells.Item(5, 17).Value = 0.08333758978133335
$ws.Cells.Item(5, 18).Value = 0.500025538688
$ws.Cells.Item(5, 19).Value = 0.003053568554213044
$ws.Cells.Item(5, 20).Value = 0.003041374773429665
$ws.Cells.Item(6, 7).Value = 0.05619466666666667
$ws.Cells.Item(6, 9).Value = 0.04986276087265156
$ws.Cells.Item(6, 10).Value = 0.07297477932340853
$ws.Cells.Item(6, 13).Value = 3.953360666666667
$ws.Cells.Item(6, 14).Value = 11.860082
$ws.Cells.Item(6, 15).Value = 0.1632495352375677
$ws.Cells.Item(6, 16).Value = 0.1666514161534747
$ws.Cells.Item(6, 17).Value = 0.2221577848764445
$ws.Cells.Item(6, 18).Value = 1.999420063888
$ws.Cells.Item(6, 19).Value = 0.008140072538122345
$ws.Cells.Item(6, 20).Value = 0.01216135031773334
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.070792
$ws.Cells.Item(7, 8).Value = 2.141584
$ws.Cells.Item(7, 9).Value = 0.9501372391273485
$ws.Cells.Item(7, 10).Value = 0.9270252206765914
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.874784666666667
$ws.Cells.Item(7, 14).Value = 8.624354
$ws.Cells.Item(7, 15).Value = 0.1187109652550681
$ws.Cells.Item(7, 16).Value = 0.121184727686443
$ws.Cells.Item(7, 17).Value = 3.078296422789333
$ws.Cells.Item(7, 18).Value = 18.469778536736
$ws.Cells.Item(7, 19).Value = 0.112791708781593
$ws.Cells.Item(7, 20).Value = 0.1123412989261575
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.070792
$ws.Cells.Item(8, 8).Value = 2.141584
$ws.Cells.Item(8, 9).Value = 0.9501372391273485
$ws.Cells.Item(8, 10).Value = 0.9270252206765914
$ws.Cells.Item(8, 15).Value = 0.4442422727481699
$ws.Cells.Item(8, 16).Value = 0.4534996302499962
$ws.Cells.Item(8, 17).Value = 11.51965529144
$ws.Cells.Item(8, 18).Value = 69.11793174863999
$ws.Cells.Item(8, 19).Value = 0.4220911265326046
$ws.Cells.Item(8, 20).Value = 0.4204055948092553
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.070792
$ws.Cells.Item(9, 8).Value = 2.141584
$ws.Cells.Item(9, 9).Value = 0.9501372391273485
$ws.Cells.Item(9, 10).Value = 0.9270252206765914
$ws.Cells.Item(9, 13).Value = 5.147441999999999
$ws.Cells.Item(9, 14).Value = 15.442326
$ws.Cells.Item(9, 15).Value = 0.2125577666737049
$ws.Cells.Item(9, 16).Value = 0.2169871588243338
$ws.Cells.Item(9, 17).Value = 5.511839714063998
$ws.Cells.Item(9, 18).Value = 33.071038284384
$ws.Cells.Item(9, 19).Value = 0.2019590495824291
$ws.Cells.Item(9, 20).Value = 0.2011525687931146
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.070792
$ws.Cells.Item(10, 8).Value = 2.141584
$ws.Cells.Item(10, 9).Value = 0.9501372391273485
$ws.Cells.Item(10, 10).Value = 0.9270252206765914
$ws.Cells.Item(10, 13).Value = 1.483016
$ws.Cells.Item(10, 14).Value = 2.966032
$ws.Cells.Item(10, 15).Value = 0.06123946008548931
$ws.Cells.Item(10, 16).Value = 0.04167706708575228
$ws.Cells.Item(10, 17).Value = 1.588001668672
$ws.Cells.Item(10, 18).Value = 6.352006674688
$ws.Cells.Item(10, 19).Value = 0.05818589153127627
$ws.Cells.Item(10, 20).Value = 0.03863569231232261
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.070792
$ws.Cells.Item(11, 8).Value = 2.141584
$ws.Cells.Item(11, 9).Value = 0.9501372391273485
$ws.Cells.Item(11, 10).Value = 0.9270252206765914
$ws.Cells.Item(11, 13).Value = 3.953360666666667
$ws.Cells.Item(11, 14).Value = 11.860082
$ws.Cells.Item(11, 15).Value = 0.1632495352375677
$ws.Cells.Item(11, 16).Value = 0.1666514161534747
$ws.Cells.Item(11, 17).Value = 4.233226974981333
$ws.Cells.Item(11, 18).Value = 25.399361849888
$ws.Cells.Item(11, 19).Value = 0.1551094626994454
$ws.Cells.Item(11, 20).Value = 0.1544900658357414
